$d = $word.ActiveDocument

# Locate the paragraph that starts the task-description block by a unique
# substring, then walk forward with .Next() so we don't depend on fixed
# paragraph indices.
$anchor = $d.Content
$found = $anchor.Find.Execute("Логическая функция от четырех переменных", `
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

$p1 = $anchor.Paragraphs(1)       # "...Протестировать ... защитить ее."

# Delete the two following paragraphs entirely (including their paragraph
# marks), which merges what follows up against $p1. Re-derive "next" from
# the still-valid $p1 each time rather than chaining stale Paragraph
# references, since those don't re-seat correctly across a mutation.
$toDelete1 = $p1.Next()            # "В соответствии с персональным вариантом ..."
$toDelete1.Range.Delete()

$toDelete2 = $p1.Next()            # "F(a, b, c, d) = D55B16"
$toDelete2.Range.Delete()

# In the remaining paragraph, change the trailing "ее." to "её.", splitting
# the tail off into its own runs (as the target XML does): the run keeps
# "...защитить е", a new run holds "ё", and another new run holds ".".
$r1 = $p1.Range
$text1 = $r1.Text
$len1 = $text1.Length
$tailStart = $r1.Start + $len1 - 3   # just after the first "е" of "ее."
$tailEnd = $r1.Start + $len1 - 1     # stop before the paragraph mark

$tail = $d.Range($tailStart, $tailEnd)
$tail.Text = ""

$yo = $d.Range($tailStart, $tailStart)
$yo.InsertAfter([string][char]0x0451)   # "ё"
$yo.Bold = $true
$yo.Bold = $false

$dotPos = $tailStart + 1
$dot = $d.Range($dotPos, $dotPos)
$dot.InsertAfter(".")
$dot.Bold = $true
$dot.Bold = $false
